$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2017045454545454
$ws.Range("C2").Value = 0.5482954545454546
$ws.Range("J2").Value = 0.009943181818181818
$ws.Range("P2").Value = 0.1534090909090909
$ws.Range("S2").Value = 0.08664772727272728
$ws.Range("B3").Value = 0.009876543209876543
$ws.Range("C3").Value = 0.03209876543209877
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.7407407407407407
$ws.Range("S3").Value = 0.1925925925925926
$ws.Range("J4").Value = 0.03896103896103896
$ws.Range("O4").Value = 0.01298701298701299
$ws.Range("P4").Value = 0.7532467532467533
$ws.Range("S4").Value = 0.1948051948051948
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.07773109243697479
$ws.Range("D6").Value = 0.01050420168067227
$ws.Range("E6").Value = 0.002100840336134454
$ws.Range("F6").Value = 0.06092436974789916
$ws.Range("J6").Value = 0.2521008403361344
$ws.Range("O6").Value = 0.01260504201680672
$ws.Range("Q6").Value = 0.1701680672268908
$ws.Range("R6").Value = 0.06512605042016807
$ws.Range("S6").Value = 0.3487394957983193
$ws.Range("B7").Value = 0.1102564102564103
$ws.Range("D7").Value = 0.02307692307692308
$ws.Range("F7").Value = 0.04871794871794872
$ws.Range("J7").Value = 0.1230769230769231
$ws.Range("O7").Value = 0.01794871794871795
$ws.Range("Q7").Value = 0.1641025641025641
$ws.Range("R7").Value = 0.05641025641025641
$ws.Range("S7").Value = 0.4564102564102564
$ws.Range("B8").Value = 0.1137123745819398
$ws.Range("D8").Value = 0.01783723522853958
$ws.Range("F8").Value = 0.05797101449275362
$ws.Range("J8").Value = 0.1270903010033445
$ws.Range("O8").Value = 0.02006688963210702
$ws.Range("Q8").Value = 0.1560758082497213
$ws.Range("R8").Value = 0.0903010033444816
$ws.Range("S8").Value = 0.4169453734671126
$ws.Range("B9").Value = 0.1119592875318066
$ws.Range("D9").Value = 0.01017811704834606
$ws.Range("F9").Value = 0.04834605597964377
$ws.Range("J9").Value = 0.1068702290076336
$ws.Range("O9").Value = 0.02290076335877863
$ws.Range("Q9").Value = 0.1806615776081425
$ws.Range("R9").Value = 0.09669211195928754
$ws.Range("S9").Value = 0.4223918575063613
$ws.Range("B10").Value = 0.1277184657967576
$ws.Range("D10").Value = 0.0166073546856465
$ws.Range("E10").Value = 0.0007908264136022143
$ws.Range("F10").Value = 0.07671016211941478
$ws.Range("J10").Value = 0.1091340450771056
$ws.Range("O10").Value = 0.02649268485567418
$ws.Range("Q10").Value = 0.2079873467773824
$ws.Range("R10").Value = 0.07236061684460261
$ws.Range("S10").Value = 0.3621984974298141
$ws.Range("G11").Value = 0.1368421052631579
$ws.Range("J11").Value = 0.1112781954887218
$ws.Range("K11").Value = 0.1909774436090226
$ws.Range("L11").Value = 0.5398496240601504
$ws.Range("S11").Value = 0.02105263157894737
$ws.Range("G12").Value = 0.7228260869565217
$ws.Range("J12").Value = 0.1847826086956522
$ws.Range("K12").Value = 0.01902173913043478
$ws.Range("L12").Value = 0.01358695652173913
$ws.Range("S12").Value = 0.05978260869565218
$ws.Range("G13").Value = 0.7027027027027027
$ws.Range("J13").Value = 0.2162162162162162
$ws.Range("S13").Value = 0.08108108108108109
$ws.Range("F15").Value = 0.01268498942917548
$ws.Range("H15").Value = 0.186046511627907
$ws.Range("I15").Value = 0.06342494714587738
$ws.Range("J15").Value = 0.3255813953488372
$ws.Range("K15").Value = 0.05496828752642706
$ws.Range("M15").Value = 0.002114164904862579
$ws.Range("O15").Value = 0.0824524312896406
$ws.Range("S15").Value = 0.2727272727272727
$ws.Range("F16").Value = 0.02188183807439825
$ws.Range("H16").Value = 0.1597374179431072
$ws.Range("I16").Value = 0.08971553610503283
$ws.Range("J16").Value = 0.3894967177242888
$ws.Range("K16").Value = 0.1094091903719912
$ws.Range("M16").Value = 0.02188183807439825
$ws.Range("O16").Value = 0.04595185995623632
$ws.Range("S16").Value = 0.161925601750547
$ws.Range("F17").Value = 0.02066590126291619
$ws.Range("H17").Value = 0.1526980482204363
$ws.Range("I17").Value = 0.08495981630309989
$ws.Range("J17").Value = 0.4064293915040184
$ws.Range("K17").Value = 0.1148105625717566
$ws.Range("M17").Value = 0.01607347876004592
$ws.Range("O17").Value = 0.05855338691159587
$ws.Range("S17").Value = 0.1458094144661309
$ws.Range("F18").Value = 0.01120448179271709
$ws.Range("H18").Value = 0.1764705882352941
$ws.Range("I18").Value = 0.09243697478991597
$ws.Range("J18").Value = 0.3809523809523809
$ws.Range("K18").Value = 0.1288515406162465
$ws.Range("M18").Value = 0.008403361344537815
$ws.Range("O18").Value = 0.05042016806722689
$ws.Range("S18").Value = 0.1512605042016807
$ws.Range("F19").Value = 0.02726281352235551
$ws.Range("H19").Value = 0.2010178117048346
$ws.Range("I19").Value = 0.07851690294438386
$ws.Range("J19").Value = 0.3511450381679389
$ws.Range("K19").Value = 0.1061432206470374
$ws.Range("M19").Value = 0.01853871319520175
$ws.Range("N19").Value = 0.0003635041802980734
$ws.Range("O19").Value = 0.06615776081424936
$ws.Range("S19").Value = 0.1508542348237005
